$d = $word.ActiveDocument
$bodyXml = @'
<w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>FICHA DE AYUDA UNIFICADA</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t>Denominacion normativa nombre ayuda</w:t></w:r></w:p><w:p><w:r><w:t>Ayuda para sufragar parcialmente los gastos de adquisición de tarjeta Bonometro para personas mayores de 65 o con discapacidad igual o superior al 64%, 2025/2026</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t>Portales</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:t>discapacidad</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:t>mayores</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t>Categoria</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:t>transporte</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t>Tipo ayuda</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:t>transporte</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t>Fecha inicio</w:t></w:r></w:p><w:p><w:r><w:t>08/07/2025</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t>Fecha fin</w:t></w:r></w:p><w:p><w:r><w:t>27/09/2025</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t>Fecha publicacion</w:t></w:r></w:p><w:p><w:r><w:t>07/07/2025</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t>Ambito territorial</w:t></w:r></w:p><w:p><w:r><w:t>Municipio de Meliana, Valencia</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t>Administracion</w:t></w:r></w:p><w:p><w:r><w:t>Ayuntamiento de Meliana (Valencia)</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t>Plazo presentacion</w:t></w:r></w:p><w:p><w:r><w:t>Hasta fecha de fin</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t>Destinatarios</w:t></w:r></w:p><w:p><w:r><w:t>Personas mayores de 65 años o con un grado de discapacidad igual o superior al 64%, empadronadas en Meliana.</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t>Requisitos acceso</w:t></w:r></w:p><w:p><w:r><w:t>1. Tener 65 años o más a fecha 31 de octubre de 2025, o certificado de minusvalía igual o superior al 64%. 2. Empadronamiento en Meliana con al menos 1 año de antigüedad a fecha 31 de octubre de 2025. 3. Estar al corriente con la Agencia Tributaria, Seguridad Social y obligaciones tributarias del Ayuntamiento de Meliana. 4. No estar incurso en las circunstancias del artículo 13 de la Ley 38/2003. 5. No incurrir en causas de exclusión o incompatibilidad especificadas en la base décima.</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t>Descripcion</w:t></w:r></w:p><w:p><w:r><w:t>La ayuda consiste en una subvención para el abono parcial del título de transporte de Ferrocarrils de la Generalitat Valenciana (FGV) para personas mayores de 65 años o con discapacidad reconocida igual o superior al 64%. La gestión y trámite del bono será realizado por el Ayuntamiento de Meliana. La ayuda es incompatible con cualquier otra subvención o ayuda, pública o privada, que tenga la misma finalidad.</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t>Cuantia</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:t>Aportación del beneficiario ................................. 50 €</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:t>Subvención para mayores de 65 años .......................... 66,40 €</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:t>Subvención para personas con discapacidad ................... 37,30 €</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:t>Gasto de reactivación de tarjeta ............................ 2 €</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:t>Gasto de expedición de tarjeta nueva ........................ 4 €</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:t>Subvención abono parcial título transporte .................. Parcial Título de transporte</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:t>Aportación personas mayores ................................. 50,00 €</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:t>Subvención personas mayores ................................. 66,40 €</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:t>Aportación personas con discapacidad ........................ 50,00 €</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:t>Subvención personas con discapacidad ........................ 37,30 €</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:t>Gasto reactivación tarjeta .................................. 2,00 €</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:t>Expedición tarjeta nueva .................................... 4,00 €</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t>Importe maximo</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:t>Concepto: Subvención máxima para mayores de 65 años • Cantidad: 66,40 €</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:t>Concepto: Subvención máxima para personas con discapacidad • Cantidad: 37,30 €</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:t>Concepto: Subvención abono parcial título transporte • Cantidad: Parcial</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:t>Concepto: Subvención máxima por persona mayor • Cantidad: 66,40 €</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:t>Concepto: Subvención máxima por persona con discapacidad • Cantidad: 37,30 €</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t>Costes no subvencionables</w:t></w:r></w:p><w:p><w:r><w:t>La percepción de esta ayuda es incompatible con cualquier otra subvención o ayuda, pública o privada, que tenga la misma finalidad.</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t>Resolucion</w:t></w:r></w:p><w:p><w:r><w:t>El órgano competente resolverá el procedimiento de concesión en un plazo mínimo de quince días naturales desde la fecha de elevación de la propuesta de resolución. El plazo de notificación de la resolución será de diez días. El silencio administrativo es negativo.</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t>Documentos presentar</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:t>Clave: Instancia-solicitud • Valor: Según modelo adjunto como anexo I.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:t>Clave: Autorización para obtener certificaciones telemáticas • Valor: AEAT y TGSS.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:t>Clave: Documento de identidad • Valor: DNI, NIE o pasaporte.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:t>Clave: Certificado de discapacidad • Valor: Igual o superior al 64%.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:t>Clave: Fotografía reciente • Valor: Formato JPG (solo para Sede Electrónica).</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:t>Clave: DNI, NIE o pasaporte • Valor: Documento de identificación personal</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t>Normativa reguladora</w:t></w:r></w:p><w:p><w:r><w:t>- Bases reguladoras y convocatoria de la ayuda para sufragar parcialmente los gastos de adquisición de tarjeta Bonometro, 2025/2026. BOP Valencia núm. 123, 07/07/2025.</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t>Referencia legislativa</w:t></w:r></w:p><w:p><w:r><w:t>- Ley 38/2003, de 17 de noviembre, General de Subvenciones.</w:t><w:br/><w:t>- Ley 38/2003, de 17 de noviembre, general de Subvenciones.</w:t><w:br/><w:t>- Ley 39/2015, de 1 de octubre, del Procedimiento Administrativo común de las Administraciones Públicas.</w:t><w:br/><w:t>- Ley 7/1985, de 2 de abril, reguladora de las Bases de Régimen Local.</w:t><w:br/><w:t>- Real Decreto 887/2006, de 21 de julio, por el que se aprueba el Reglamento de la Ley General de Subvenciones.</w:t><w:br/><w:t>- Real Decreto Legislativo 2/2004, de 5 de marzo, por el que se aprueba el Texto refundido de la Ley Reguladora de las Haciendas Locales.</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t>Lugares presentacion</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t>Presencial</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:t>Clave: Presencialmente en: • Valor: Centro de Servicios Sociales de Atención Primaria del Ayuntamiento de Meliana, calle Ramón y Cajal 4, 46133 Meliana (Valencia). - Oficinas de correos para Registro de documentos según procedimiento administrativo (con el sobre abierto para compulsa de documentos). - Registro de Ventanilla Única en el territorio Nacional.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:t>Clave: Presencialmente en: • Valor: Ayuntamiento de Meliana, Dirección: Plaza Mayor, 1, 46133 Meliana, Valencia. Teléfono: 961 490 000. Oficinas de correos para Registro de documentos según procedimiento administrativo (con el sobre abierto para compulsa de documentos). Registro de Ventanilla Única en el territorio Nacional.</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t>Online</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:t>Clave: Electrónicamente en: • Valor: https://meliana.sede.dival.es/. También puede presentarse a través de la Red SARA (Sistema de Aplicaciones y Redes para las Administraciones), una plataforma estatal segura que permite enviar solicitudes electrónicas a cualquier administración pública, usando certificado digital o sistema Cl@ve.</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t>Criterios concesion</w:t></w:r></w:p><w:p><w:r><w:t>La concesión se realizará mediante el procedimiento de concurrencia competitiva.</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t>Usuario</w:t></w:r></w:p><w:p><w:r><w:t>MELIANA</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t>Fecha</w:t></w:r></w:p><w:p><w:r><w:t>07/07/2025</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t>Frase publicitaria</w:t></w:r></w:p><w:p><w:r><w:t>Solicita la Tarjeta Bonometro 2025/2026 y ahorra en transporte si eres mayor de 65 o tienes discapacidad. Infórmate ya.</w:t></w:r></w:p>
'@
$d.Content.InsertXML($bodyXml)
Write-Output "done"
